# Office supplies acquisition table.xlsx — apply authored edits:
#  - rename the (only) worksheet from "Sheet1" to "1.2001"
#  - move the active selection on that sheet from I13 to J32

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "1.2001"

# Move/restore the selection to cell J32 (also updates ActiveCell).
$ws.Range("J32").Select()
